$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 89 is a new data row identical to row 88 except for the date (A89).
# Copy the existing row 88 (which already holds all values as text) down to
# row 89, then just patch the date so all of the non-date columns keep
# their original text formatting untouched.
$ws.Range("A88:J88").Copy()
$ws.Range("A89").PasteSpecial()

# Force the date cell to remain plain text (matching the source file, where
# every value - including the date - is stored as text) and update the day
# from 28 to 29.
$ws.Range("A89").NumberFormat = "@"
$ws.Range("A89").Characters(10, 1).Text = "9"

$excel.CutCopyMode = $false
